# RPA datasets push 2024-07-13
# Insert a new data row (LS / 이베스트스팩6호) at row 6, shifting the
# existing rows 6-27 down to 7-28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 6:27 down to make room for the new record.
$ws.Rows(6).Insert()

# Force text formatting on the date-like cells before assigning, so Excel's
# automatic date recognition doesn't turn "2024-07-02" / "2024-07-05" /
# "2024-07-12" into date serial numbers, then reset the style back to the
# workbook's default ("Normal") so the new row matches the plain, unstyled
# data cells used throughout the rest of the sheet.
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "2024-07-02"
$ws.Range("B6").Style = "Normal"

$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "2024-07-05"
$ws.Range("F6").Style = "Normal"

$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "2024-07-12"
$ws.Range("G6").Style = "Normal"

$ws.Range("A6").Value = "LS"
$ws.Range("C6").Value = "이베스트스팩6호"
$ws.Range("D6").Value = "엘에스"
$ws.Range("E6").Value = "엘에스"

$ws.Range("H6").Value = 8000
$ws.Range("I6").Value = 4000000
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 100
